$d = $word.ActiveDocument

# --- Re-balance the FSM state-transition table's column widths --------------
# (old widths, in dxa) 1004 1005 581 1005 1005 1005 1005 1005 1005 724
# (new widths, in dxa)  998  997 606  997  997  997  997  997  997 761
$t = $d.Tables.Item(1)
$newWidthsDxa = @(998, 997, 606, 997, 997, 997, 997, 997, 997, 761)
for ($i = 1; $i -le $newWidthsDxa.Length; $i++) {
    $col = $t.Columns.Item($i)
    $col.Width = $newWidthsDxa[$i - 1] / 20.0
}

# --- Drop the stale "Idle; " prefix from the state-cell labels --------------
$d.Content.Find.Execute("Idle; Show Unselected", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Show Unselected", 2)
$d.Content.Find.Execute("Idle; Show Selected", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Show Selected", 2)
